# Applies the hermes.xlsx "Add files via upload" edit: a handful of rows in
# the phone-signup log were resynced (new entries inserted/reordered) and one
# brand-new row was appended at the end, growing the sheet from 195 to 196
# data+header rows (A1:C195 -> A1:C196). Only the affected rows are touched;
# everything else is left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phone numbers / DDD codes / dates are stored as plain TEXT in the sheet
# (e.g. "+5511967663538", "11", "2024-10-15"). Force the number format to
# "@" (Text) on each touched row before writing, otherwise Excel would
# auto-coerce these numeric-looking strings into real numbers/dates and
# silently drop things like the leading "+".

$ws.Range("A2:C9").NumberFormat = "@"
$ws.Range("A19:C22").NumberFormat = "@"
$ws.Range("A187:C196").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "+5511967663538"
$ws.Cells.Item(2, 2).Value = "11"
$ws.Cells.Item(2, 3).Value = "2024-10-15"

$ws.Cells.Item(3, 1).Value = "+5516999928423"
$ws.Cells.Item(3, 2).Value = "16"
$ws.Cells.Item(3, 3).Value = "2024-10-15"

$ws.Cells.Item(4, 1).Value = "+5511975292030"
$ws.Cells.Item(4, 2).Value = "11"
$ws.Cells.Item(4, 3).Value = "2024-10-14"

$ws.Cells.Item(5, 1).Value = "+553491814035"
$ws.Cells.Item(5, 2).Value = "34"
$ws.Cells.Item(5, 3).Value = "2024-10-14"

$ws.Cells.Item(6, 1).Value = "+556392880287"
$ws.Cells.Item(6, 2).Value = "63"
$ws.Cells.Item(6, 3).Value = "2024-10-10"

$ws.Cells.Item(7, 1).Value = "+5516996469888"
$ws.Cells.Item(7, 2).Value = "16"
$ws.Cells.Item(7, 3).Value = "2024-10-10"

$ws.Cells.Item(8, 1).Value = "+555199199744"
$ws.Cells.Item(8, 2).Value = "51"
$ws.Cells.Item(8, 3).Value = "2024-10-10"

$ws.Cells.Item(9, 1).Value = "+5511958531213"
$ws.Cells.Item(9, 2).Value = "11"
$ws.Cells.Item(9, 3).Value = "2024-10-09"

$ws.Cells.Item(19, 1).Value = "+5511974585735"
$ws.Cells.Item(19, 2).Value = "11"
$ws.Cells.Item(19, 3).Value = "2024-10-08"

$ws.Cells.Item(20, 1).Value = "+5516991966214"
$ws.Cells.Item(20, 2).Value = "16"
$ws.Cells.Item(20, 3).Value = "2024-10-08"

$ws.Cells.Item(21, 1).Value = "+5511981050835"
$ws.Cells.Item(21, 2).Value = "11"
$ws.Cells.Item(21, 3).Value = "2024-10-08"

$ws.Cells.Item(22, 1).Value = "+5511981996340"
$ws.Cells.Item(22, 2).Value = "11"
$ws.Cells.Item(22, 3).Value = "2024-10-07"

$ws.Cells.Item(187, 1).Value = "+5521985109311"
$ws.Cells.Item(187, 2).Value = "21"
$ws.Cells.Item(187, 3).Value = "2024-09-09"

$ws.Cells.Item(188, 1).Value = "+5511957562684"
$ws.Cells.Item(188, 2).Value = "11"
$ws.Cells.Item(188, 3).Value = "2024-09-09"

$ws.Cells.Item(189, 1).Value = "+5524999327754"
$ws.Cells.Item(189, 2).Value = "24"
$ws.Cells.Item(189, 3).Value = "2024-09-09"

$ws.Cells.Item(190, 1).Value = "+558699687586"
$ws.Cells.Item(190, 2).Value = "86"
$ws.Cells.Item(190, 3).Value = "2024-09-09"

$ws.Cells.Item(191, 1).Value = "+556185975181"
$ws.Cells.Item(191, 2).Value = "61"
$ws.Cells.Item(191, 3).Value = "2024-08-30"

$ws.Cells.Item(192, 1).Value = "+553291004823"
$ws.Cells.Item(192, 2).Value = "32"
$ws.Cells.Item(192, 3).Value = "2024-08-26"

$ws.Cells.Item(193, 1).Value = "+34603138909"
# This entry has no DDD code in the source data (B column is blank there);
# a lone quote-prefix yields an empty TEXT cell instead of an empty/General one.
$ws.Cells.Item(193, 2).Value = "'"
$ws.Cells.Item(193, 3).Value = "2024-08-19"

$ws.Cells.Item(194, 1).Value = "+5511981274889"
$ws.Cells.Item(194, 2).Value = "11"
$ws.Cells.Item(194, 3).Value = "2024-08-01"

$ws.Cells.Item(195, 1).Value = "+5511967859426"
$ws.Cells.Item(195, 2).Value = "11"
$ws.Cells.Item(195, 3).Value = "2024-07-28"

$ws.Cells.Item(196, 1).Value = "+5521965197022"
$ws.Cells.Item(196, 2).Value = "21"
$ws.Cells.Item(196, 3).Value = "2024-07-21"
